# Applies two text edits to the document:
#  1. Paragraph "... onde conste o total de horas trabalhadas para cada func|ionário."
#     - the run split (caused by the stray "_GoBack" bookmark sitting between
#       "func" and "ionário") is removed and the two runs are merged back
#       into a single run.
#  2. Paragraph "O acesso ao sistema ... (Administrador, RH, Gerente, Funcionário), ..."
#     - ", Funcionário" is removed from the parenthetical list, and the
#       sentence is re-split into three runs with the (now relocated)
#       "_GoBack" bookmark sitting between "...Gerente" and "), onde ...".

$d = $word.ActiveDocument

# ============================================================
# Part 1: paragraph 10 - merge the "func" / "ionário." runs and
# drop the bookmark that used to sit between them.
# ============================================================

$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$full = $d.Content.Text
$consteIdx = $full.IndexOf("conste")
$afterConsteIdx = $consteIdx + "conste".Length

# Temporarily bookmark the boundary right after "conste" so that the
# upcoming same-text "replace" (which coalesces contiguous, identically
# formatted runs) stops there instead of also swallowing the "conste" run.
$d.Bookmarks.Add("TEMPBOUNDARY", $d.Range($afterConsteIdx, $afterConsteIdx))

$phrase10 = " o total de horas trabalhadas para cada funcionário."
$d.Content.Find.Execute($phrase10, $true, $false, $false, $false, $false, $true, 1, $false, $phrase10, 2) | Out-Null

$tb = $d.Bookmarks("TEMPBOUNDARY")
$tb.Delete()

# ============================================================
# Part 2: paragraph 8 - drop ", Funcionário" from the list and
# re-split the sentence into three runs around the new bookmark spot.
# ============================================================

$d.Content.Find.Execute(", Funcionário)", $true, $false, $false, $false, $false, $true, 1, $false, ")", 2) | Out-Null

$full2 = $d.Content.Text
$paraStart = $full2.IndexOf("O acesso ao sistema")
$s1 = "O acesso ao sistema pelos usuários deverá ser condicionado a níveis de acesso (Adminis"
$s2 = "trador, RH, Gerente"
$split1 = $paraStart + $s1.Length
$split2 = $split1 + $s2.Length

# Toggling Bold on and back off over the middle segment forces Word to
# split the single run into three runs at $split1/$split2 without leaving
# any residual formatting difference behind.
$middle = $d.Range($split1, $split2)
$middle.Bold = 1
$middle.Bold = 0

# Re-create "_GoBack" at its new location, between the second and third run.
$d.Bookmarks.Add("_GoBack", $d.Range($split2, $split2))
